$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "'250.06"
Set-TextValue $ws.Range("D3") "'22.95"
Set-TextValue $ws.Range("D4") "'5.434"
Set-TextValue $ws.Range("D5") "'0.05620"
Set-TextValue $ws.Range("D6") "'3.419"
Set-TextValue $ws.Range("D7") "'6.373"
Set-TextValue $ws.Range("D8") "'0.8166"
Set-TextValue $ws.Range("D9") "'0.9142"
Set-TextValue $ws.Range("D10") "'0.1439"
Set-TextValue $ws.Range("D11") "'0.07544"
Set-TextValue $ws.Range("D12") "'0.03131"
Set-TextValue $ws.Range("D13") "'0.03087"
Set-TextValue $ws.Range("D14") "'0.09317"
Set-TextValue $ws.Range("D15") "'3.565"
Set-TextValue $ws.Range("D16") "'0.001633"
Set-TextValue $ws.Range("D17") "'0.04758"
Set-TextValue $ws.Range("B18") "'One"
Set-TextValue $ws.Range("C18") "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D18") "'0.0005801"
Set-TextValue $ws.Range("E18") "'17OneONEWorstin24h"
Set-TextValue $ws.Range("B19") "'TigerCash"
Set-TextValue $ws.Range("C19") "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D19") "'0.006424"
Set-TextValue $ws.Range("E19") "'18TigerCashTCH"
Set-TextValue $ws.Range("B20") "'HotbitToken"
Set-TextValue $ws.Range("C20") "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws.Range("D20") "'0.004989"
Set-TextValue $ws.Range("E20") "'19HotbitTokenHTB"
Set-TextValue $ws.Range("B21") "'BitKan"
Set-TextValue $ws.Range("C21") "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Range("D21") "'0.001036"
Set-TextValue $ws.Range("E21") "'20BitKanKAN"
Set-TextValue $ws.Range("B22") "'NitroEx"
Set-TextValue $ws.Range("C22") "'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws.Range("D22") "'0.0001504"
Set-TextValue $ws.Range("E22") "'21NitroExNTX"
Set-TextValue $ws.Range("B23") "'LEO"
Set-TextValue $ws.Range("C23") "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D23") "'3.722"
Set-TextValue $ws.Range("E23") "'22LEOLEO"
Set-TextValue $ws.Range("B24") "'BTSEToken"
Set-TextValue $ws.Range("C24") "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D24") "'2.187"
Set-TextValue $ws.Range("E24") "'23BTSETokenBTSE"
Set-TextValue $ws.Range("D25") "'0.3303"
Set-TextValue $ws.Range("D26") "'0.1274"
Set-TextValue $ws.Range("E27") "'26AAXTokenAAB"
Set-TextValue $ws.Range("D28") "'0.0003010"
Set-TextValue $ws.Range("D40") "'0.04016"
Set-TextValue $ws.Range("B41") "'BKEXToken"
Set-TextValue $ws.Range("C41") "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D41") "'0.1067"
Set-TextValue $ws.Range("E41") "'40BKEXTokenBKK"
Set-TextValue $ws.Range("B42") "'CEJI"
Set-TextValue $ws.Range("C42") "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D42") "'0.002728"
Set-TextValue $ws.Range("E42") "'41CEJICEJI"
Set-TextValue $ws.Range("B43") "'KickToken"
Set-TextValue $ws.Range("C43") "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D43") "'0.002944"
Set-TextValue $ws.Range("E43") "'42KickTokenKICK"
Set-TextValue $ws.Range("D45") "'0.00005582"
Set-TextValue $ws.Range("D48") "'0.5008"
Set-TextValue $ws.Range("D49") "'0.2331"
Set-TextValue $ws.Range("E49") "'48BOLOBOLOBestin24h"
Set-TextValue $ws.Range("D51") "'0.01012"
